$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the old "Terms Typically Offered" column (D),
# shifting it to G, and leaving D:F empty for the new columns
# (Corequisites, Concurrent, Recommended).
$ws.Range("D1:F1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Row 2 - UNIV 100
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"

# Row 3 - UNIV 125
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"

# Row 4 - UNIV 321
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"

# Row 5 - UNIV 330
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"

# Row 6 - UNIV 333
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "NA"

# Row 7 - UNIV 350
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "NA"

# Row 8 - UNIV 391
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"

# Row 9 - UNIV 392: the "Recommended" clause moves out of Prerequisites (C)
# into the new Recommended column (F); Terms Typically Offered (G) gains a
# trailing space to match the source data.
$ws.Range("C9").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of GE Area B1 with a grade of C- or better in at least one of the courses; and completion of GE Areas B2, B3, and B4."
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "UNIV 391 and completion of GE Areas D2 and D3."
$ws.Range("G9").Value = "SP "

# Row 10 - UNIV 424
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"

# Row 11 - UNIV 470
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "NA"

# Row 12 - UNIV 491: the "Corequisite" clause moves out of Prerequisites (C)
# into the new Corequisites column (D); Terms Typically Offered (G) gains a
# trailing space to match the source data.
$ws.Range("C12").Value = "Consent of instructor, and senior or graduate standing."
$ws.Range("D12").Value = "GE Area D5."
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "F "

# Row 13 - UNIV 492: the "Recommended" clause moves out of Prerequisites (C)
# into the new Recommended column (F); Terms Typically Offered (G) gains a
# trailing space to match the source data.
$ws.Range("C13").Value = "Junior standing and completion of GE Area B, or graduate standing."
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "UNIV 391, GE Area D2, and GE Area D3."
$ws.Range("G13").Value = "SP "
